$d = $word.ActiveDocument

$replacements = @(
    @("68÷4=17, 0", "80÷8=10, 0"),
    @("64÷5=12, 4", "68÷6=11, 2"),
    @("49÷6=8, 1", "12÷5=2, 2"),
    @("88÷5=17, 3", "64÷7=9, 1"),
    @("30÷7=4, 2", "76÷2=38, 0"),
    @("32÷4=8, 0", "17÷7=2, 3"),
    @("69÷4=17, 1", "79÷6=13, 1"),
    @("40÷6=6, 4", "39÷3=13, 0"),
    @("38÷8=4, 6", "91÷7=13, 0"),
    @("64÷3=21, 1", "91÷2=45, 1"),
    @("98÷8=12, 2", "86÷9=9, 5"),
    @("75÷4=18, 3", "25÷8=3, 1"),
    @("77÷7=11, 0", "12÷8=1, 4"),
    @("27÷6=4, 3", "38÷5=7, 3"),
    @("93÷2=46, 1", "60÷8=7, 4"),
    @("51÷5=10, 1", "99÷5=19, 4"),
    @("92÷6=15, 2", "30÷4=7, 2"),
    @("93÷6=15, 3", "85÷6=14, 1"),
    @("43÷7=6, 1", "98÷6=16, 2"),
    @("85÷8=10, 5", "15÷5=3, 0"),
    @("35÷8=4, 3", "87÷8=10, 7"),
    @("88÷8=11, 0", "69÷7=9, 6"),
    @("89÷2=44, 1", "72÷8=9, 0"),
    @("74÷4=18, 2", "56÷6=9, 2"),
    @("67÷5=13, 2", "81÷4=20, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
